$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the paragraph ending in "=embed)" to the
#    very first (empty) paragraph of the document. Word only allows a single
#    bookmark with a given name in a document, so re-adding "_GoBack" at the
#    new location automatically removes it from its old location.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$d.Bookmarks.Add("_GoBack", $firstPara.Range)

# ---------------------------------------------------------------------------
# 2) Split the run containing "embed.set_footer" into "embed.set_" + "footer"
#    (formatting stays identical on both halves). We force Word to materialize
#    two separate runs by toggling a character property on the "footer" part
#    and then clearing it again, which leaves the visible formatting
#    untouched but breaks the run in two.
# ---------------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("embed.set_footer")
if ($found) {
    $callStart = $find.Start
    $callEnd = $find.End
    $splitPoint = $callStart + "embed.set_".Length
    $footerPart = $d.Range($splitPoint, $callEnd)
    $footerPart.Bold = 1
    $footerPart.Bold = 0
}

# ---------------------------------------------------------------------------
# 3) Replace the footer text literal with the new credits message. The new
#    run needs to end up with *only* a color override (no rFonts/lang), which
#    is what happens naturally when it is typed fresh and never touches
#    Font.Name/Language setters. We build it in a throwaway paragraph first
#    (so it starts out with no direct character formatting at all besides the
#    color we set), copy it across via FormattedText (which carries the
#    run-level formatting, unlike plain .Text), and then remove the scratch
#    paragraph again.
# ---------------------------------------------------------------------------
$newCreditsText = "Check out the ``credits`` command to see the list of people who helped making this bot."

$paraCountBefore = $d.Paragraphs.Count
$d.Paragraphs.Add() | Out-Null
$scratchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratchPara.Range.InsertAfter($newCreditsText)
$scratchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$cleanRange = $d.Range($scratchPara.Range.Start, $scratchPara.Range.Start + $newCreditsText.Length)
$cleanRange.Font.Color = 5867370   # 0x6A8759 packed as 0x00BBGGRR

$body = $d.Content
$oldCredits = '"Made by Mando_The_Mercenary#9484 | Bot''s Logo designed by ARCAS#0954 | Build info provided by SIGMA#5422"'
$bodyFound = $body.Find.Execute($oldCredits)
if ($bodyFound) {
    $quoteStart = $body.Start
    $quoteEnd = $body.End
    $innerRange = $d.Range($quoteStart + 1, $quoteEnd - 1)
    $innerRange.FormattedText = $cleanRange.FormattedText
}

$paraCountNow = $d.Paragraphs.Count
$prevEnd = $d.Paragraphs.Item($paraCountNow - 1).Range.End
$scratchEnd = $d.Paragraphs.Item($paraCountNow).Range.End
$cleanupRange = $d.Range($prevEnd, $scratchEnd)
$cleanupRange.Delete()

Write-Output "paragraphs before=$paraCountBefore after=$($d.Paragraphs.Count)"
